# Generate Report for Handback
# The handback report swaps which source file ("a7ba330e-...md" vs
# "84b09259-...md") is listed in row 2 vs row 3 of each sheet, marks the
# 84b09259 entry as "Handed back: in sync with en-US" (it previously said
# "Ready for handoff"), and records fresh "Latest Handback DateTime"
# timestamps for the newly handed-back file in the locale sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.md"
$ov.Range("A3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.md"
$zh.Range("C2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.ebc2fbcc84a15bf51808aaae1467fa9b88a4ceb3.zh-cn.xlf"
$zh.Range("D2").Value = "2016-01-26 09:06:57"
$zh.Range("E2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.md"
$zh.Range("F2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.ebc2fbcc84a15bf51808aaae1467fa9b88a4ceb3.zh-cn.xlf"
$zh.Range("G2").Value = "2016-01-26 09:07:41"

$zh.Range("A3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.fe500b8130d0f1a52ea71db28b4502d3da31c4d8.zh-cn.xlf"
$zh.Range("D3").Value = "2016-01-26 09:05:13"
$zh.Range("E3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md"
$zh.Range("F3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.fe500b8130d0f1a52ea71db28b4502d3da31c4d8.zh-cn.xlf"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.md"
$de.Range("C2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.ebc2fbcc84a15bf51808aaae1467fa9b88a4ceb3.de-de.xlf"
$de.Range("D2").Value = "2016-01-26 09:07:08"
$de.Range("E2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.md"
$de.Range("F2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.ebc2fbcc84a15bf51808aaae1467fa9b88a4ceb3.de-de.xlf"
$de.Range("G2").Value = "2016-01-26 09:07:59"

$de.Range("A3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md"
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.fe500b8130d0f1a52ea71db28b4502d3da31c4d8.de-de.xlf"
$de.Range("D3").Value = "2016-01-26 09:05:24"
$de.Range("E3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md"
$de.Range("F3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.fe500b8130d0f1a52ea71db28b4502d3da31c4d8.de-de.xlf"
